$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for rows 2-41.
# Update every row where the existing value is 45185 (2023-09-16) to 45204 (2023-10-05).
for ($row = 2; $row -le 41; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value = 45204
    }
}
